$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 4 (old MuSCs-sourced row) entirely - table shrinks from 4 data rows to 2
$ws.Rows("4:4").Delete()

# Row 2 (FAPs -> Ccl12/Ccr3 -> Resolving-Mac) keeps its text labels, only the
# numeric NATMI statistics were recomputed with the new TPM values.
$ws.Range("I2").Value = 0.02394963654761903
$ws.Range("J2").Value = 0.02394963654761903
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3389413333333333
$ws.Range("N2").Value = 1.016824
$ws.Range("Q2").Value = 0.290522885984
$ws.Range("R2").Value = 2.614705973856
$ws.Range("S2").Value = 0.02394963654761903
$ws.Range("T2").Value = 0.02394963654761903

# Row 3 used to be the MuSCs sending-cluster row; that cluster was dropped
# from the analysis entirely, so row 3 now carries the (recomputed)
# Resolving-Mac -> Ccl12/Ccr3 -> Resolving-Mac autocrine edge that used to
# live in row 4.
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 34.93245566666667
$ws.Range("H3").Value = 104.797367
$ws.Range("I3").Value = 0.9760503634523809
$ws.Range("J3").Value = 0.9760503634523809
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3389413333333333
$ws.Range("N3").Value = 1.016824
$ws.Range("Q3").Value = 11.84005310026756
$ws.Range("R3").Value = 106.560477902408
$ws.Range("S3").Value = 0.9760503634523809
$ws.Range("T3").Value = 0.9760503634523809
